$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values (keep headers in row 4 and C5/D5/G5 unchanged)
$ws.Range("E5").Value = 6
$ws.Range("H5").Value = "Box"
$ws.Range("I5").Value = "Clip"
$ws.Range("F5").Value = "Clips Paper Small"

# Delete row 6 entirely, shifting cells up
$ws.Rows("6:6").Delete()

# Update the selected range to match the target view state
$ws.Range("A5:XFD5").Select()
